$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark from the title paragraph ("YIFAN LI").
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Re-add the "_GoBack" bookmark (collapsed, zero-length) to the empty
#    paragraph that immediately precedes the first table (the empty line
#    right before the education/skills table).
$p16 = $d.Paragraphs.Item(16)
$d.Bookmarks.Add("_GoBack", $p16.Range)

# 3. Skills line: "Eclipse, Git, Gerrit, Maven, Postgres, etc." -> remove
#    "Gerrit" and the trailing "etc."
$d.Content.Find.Execute(", Gerrit, Maven,", $true, $false, $false, $false, `
                         $false, $true, 1, $false, ", Maven,", 2) | Out-Null

$d.Content.Find.Execute(", etc.", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "", 2) | Out-Null

# 4. Operating systems line: drop the trailing period after "Mac OS X".
$d.Content.Find.Execute(", Windows, Mac OS X.", $true, $false, $false, $false, `
                         $false, $true, 1, $false, ", Windows, Mac OS X", 2) | Out-Null
